$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column (A) for the "TabName" labels, shifting the
# existing query/dbExcel/WebExcel columns one place to the right (B:E).
$ws.Columns("A").Insert()

# --- Row 1 header values ---
$ws.Range("A1").Value = "TabName"
# B1 ("query"), C1 ("StatQuery"), D1 ("dbExcel") and E1 ("WebExcel") already
# have the correct values after the column insert shifted them over.

# --- Long Neo4j query text used for the Case/Sample/File tabs ---
$qCaseId = @'
MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE s.clinical_study_designation IN ['NCATS-COP01'] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'') AS `Case ID` , coalesce(s.clinical_study_designation,'') AS `Study Code` , coalesce(s.clinical_study_type,'') AS  `Study Type`, coalesce(demo.breed,'') AS Breed , coalesce(diag.disease_term,'') AS Diagnosis , coalesce(diag.stage_of_disease,'') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'') AS Age , coalesce(demo.sex,'') AS Sex , coalesce(demo.neutered_indicator,'') AS  `Neutered Status`
'@

$qStat = @'
 MATCH (p:program)<--(s:study)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
      WHERE (size([]) = 0 OR s.clinical_study_designation IN [])
        AND (s.study_disposition = 'Unrestricted')
        AND (size([]) = 0 OR s.clinical_study_type IN [])
        AND (size(['NCATS-COP01']) = 0 OR demo.breed IN ['NCATS-COP01'])
        AND (size([]) = 0 OR demo.sex IN [])
        AND (size([]) = 0 OR demo.neutered_indicator IN [])
        AND (size([]) = 0 OR diag.disease_term IN [])
        AND (size([]) = 0 OR diag.primary_disease_site IN [])
        AND (size([]) = 0 OR diag.stage_of_disease IN [])
        AND (size([]) = 0 OR diag.best_response IN [])
    OPTIONAL MATCH (c)-->(co:cohort)
    OPTIONAL MATCH (f:file)-[*]->(c)
    OPTIONAL MATCH (f)-->(parent)
    OPTIONAL MATCH (samp:sample)-->(c)
    OPTIONAL MATCH (samp)<--(al:aliquot)
    WITH DISTINCT c AS c, p, s, co, demo, diag, f, parent, samp, al
      WHERE (size([]) = 0 OR samp.summarized_sample_type IN [])
        AND (size([]) = 0 OR samp.specific_sample_pathology IN [])
        AND (size([]) = 0 OR samp.sample_site IN [])
        AND (size([]) = 0 OR head(labels(parent)) IN [])
        AND (size([]) = 0 OR f.file_type IN [])
        AND (size([]) = 0 OR f.file_format IN [])
    WITH c.case_id AS case_id,
         s.clinical_study_designation AS study_code,
         s.clinical_study_type AS study_type,
         co.cohort_description AS cohort,
         demo.breed AS breed,
         diag.disease_term AS diagnosis,
         diag.stage_of_disease AS stage_of_disease,
         diag.primary_disease_site AS disease_site,
         demo.patient_age_at_enrollment AS age,
         demo.sex AS sex,
         demo.neutered_indicator AS neutered_status,
         demo.weight AS weight,
         diag.best_response AS response_to_treatment,
         samp.sample_id AS sample_id,
         f.uuid AS file_id,
         al
    RETURN
COUNT(DISTINCT file_id) as number_of_files,
COUNT(DISTINCT sample_id) as number_of_sample,
COUNT(DISTINCT case_id) as number_of_cases,
COUNT(DISTINCT study_code) as number_of_study,
COUNT(DISTINCT al) as number_of_aliquot
    
'@

# --- Row 2: CasesTab ---
$ws.Range("A2").Value = "CasesTab"
$ws.Range("B2").Style = $ws.Range("A1").Style
$ws.Range("C2").Style = $ws.Range("A1").Style
$ws.Range("B2").Value = $qCaseId
$ws.Range("C2").Value = $qStat
$ws.Range("D2").Value = "TC03_Canine_Filter_Study-GLIOMA_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC03_Canine_Filter_Study-GLIOMA_WebData.xlsx"

# --- Row 3: SamplesTab ---
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("B3").Style = $ws.Range("A1").Style
$ws.Range("C3").Style = $ws.Range("A1").Style
$ws.Range("B3").Value = $qCaseId
$ws.Range("C3").Value = $qStat
$ws.Range("D3").Value = "TC03_Canine_Filter_Study-GLIOMA_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC03_Canine_Filter_Study-GLIOMA_WebData.xlsx"

# --- Row 4: FilesTab ---
$ws.Range("A4").Value = "FilesTab"
$ws.Range("B4").Style = $ws.Range("A1").Style
$ws.Range("C4").Style = $ws.Range("A1").Style
$ws.Range("B4").Value = $qCaseId
$ws.Range("C4").Value = $qStat
$ws.Range("D4").Value = "TC03_Canine_Filter_Study-GLIOMA_Neo4jData.xlsx"
$ws.Range("E4").Value = "TC03_Canine_Filter_Study-GLIOMA_WebData.xlsx"

# Re-apply the wrap-text cell style (same "Normal 2" style used originally)
# to the new query cells so the long text wraps inside the cell.
$ws.Range("B2:C4").WrapText = $true

# The query cells now hold very large blocks of text; grow the rows to
# Excel's maximum row height (as happens automatically once the text is
# long enough to need it), matching the source workbook.
$ws.Rows("2:4").RowHeight = 409.6

# Column A is new and narrow (short "XxxTab" labels) - size it to fit them.
$ws.Columns("A").ColumnWidth = 10
# Column E ("WebExcel" / *_WebData.xlsx) narrows slightly to fit its new text.
$ws.Columns("E").ColumnWidth = 37.2

# Restore the view: scrolled down to the last row, with cell B4 selected,
# and zoomed in to 115% (as in the edited workbook).
$ws.Range("A4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 115
$ws.Range("B4").Select()

Write-Host "Edit applied"
